# Applies the "Add files via upload" update to REVER_DailyTracker_BALRAJ.xlsx
# - Targets the "NOV-2020" sheet (3rd tab)
# - Fills in row 30 (previously blank placeholder row) with task #17 data
# - Adds 4 new sub-bullet rows (31-34) with Invoice/warranty/RPA related notes
# - Adjusts the sheet view (scroll position / selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOV-2020")
$ws.Activate()

$xlPasteFormats = -4122

# ---- Row 30: new task entry #17 ----
$ws.Range("A30").Value = 17
$ws.Range("B30").Value = 44148
$ws.Range("C30").Value = "RPA RLOGIC"
$ws.Range("D30").Value = "1. Downloaded DRS and Warranty Files from GSPN to generate Calls count report etc for October-2020"
$ws.Range("E30").Value = 1
$ws.Range("F30").Value = "Completed"

# B30 needs the existing date-formatted / bordered style (like B2); E30 needs the
# existing percent-formatted / bordered style (like E3). Copy formats so the
# existing style indexes get reused instead of minting duplicate styles.
$ws.Range("B2").Copy()
$ws.Range("B30").PasteSpecial($xlPasteFormats)
$ws.Range("E3").Copy()
$ws.Range("E30").PasteSpecial($xlPasteFormats)

# ---- Row 31: sub-bullet 2 ----
$ws.Range("D31").Value = "2. Downloaded InvoiceUpdatePDF to generate InvoiceUpdate summary report with client provided data"
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = "Completed"
$ws.Range("E31").NumberFormat = "0%"
$ws.Range("F10").Copy()
$ws.Range("F31").PasteSpecial($xlPasteFormats)

# ---- Row 32: sub-bullet 3 ----
$ws.Range("D32").Value = "3.Extended warranty report generated from the GSPN data"
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = "Completed"
$ws.Range("E31").Copy()
$ws.Range("E32").PasteSpecial($xlPasteFormats)
$ws.Range("F10").Copy()
$ws.Range("F32").PasteSpecial($xlPasteFormats)

# ---- Row 33: sub-bullet 4 ----
$ws.Range("D33").Value = "4. RPA Management Model report generated"
$ws.Range("E33").Value = 1
$ws.Range("F33").Value = "Completed"
$ws.Range("E31").Copy()
$ws.Range("E33").PasteSpecial($xlPasteFormats)
$ws.Range("F10").Copy()
$ws.Range("F33").PasteSpecial($xlPasteFormats)

# ---- Row 34: sub-bullet 5 ----
$ws.Range("D34").Value = "5.Tested with Mr Rahaman to verify the data"
$ws.Range("E34").Value = 1
$ws.Range("F34").Value = "Completed"
$ws.Range("E31").Copy()
$ws.Range("E34").PasteSpecial($xlPasteFormats)
$ws.Range("F10").Copy()
$ws.Range("F34").PasteSpecial($xlPasteFormats)

# ---- Row height tweaks on rows 20 and 25 (wrap-computed heights shrank) ----
$ws.Rows("20").RowHeight = 28.8
$ws.Rows("25").RowHeight = 28.8

# ---- Update the view: scroll position and active selection ----
$ws.Range("D35").Select()
$excel.ActiveWindow.ScrollRow = 22
